# The "Förändrad" (Changed) column C, for every data row (rows 2-238),
# is being bumped from serial date 45188 (2023-09-19) to 45189 (2023-09-20).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C238")
$range.Value = 45189
